# feat: add custom scenario for transportation
#
# The "GWPb" (E) and "GWP-LULUC" (F) header columns already exist but
# have no data for the transport rows. Populate them with 0 for each
# transport mode (rows 2-5), and remove the now-unused trailing
# "Notes" column (K), which has no data either.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "GWPb" / "GWP-LULUC" data cells.
$ws.Range("E2:F5").Value = 0

# Remove the trailing "Notes" column entirely.
$ws.Range("K:K").Delete() | Out-Null

# Update the active selection to match the new layout.
$ws.Range("A2").Select() | Out-Null
